# Update ORN Yearly Financials worksheet with latest data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ORN")

# --- Inventory (row 44): new values replacing prior year figures ---
$ws.Range("D44").Value = 50400
$ws.Range("E44").Value = 45400
$ws.Range("F44").Value = 64500
$ws.Range("G44").Value = 51100
$ws.Range("H44").Value = 28400
$ws.Range("I44").Value = 23600
$ws.Range("J44").Value = 18500

# --- Other Current Assets (row 45): new values replacing prior year figures ---
$ws.Range("D45").Value = 4100
$ws.Range("E45").Value = 12300
$ws.Range("F45").Value = 14100
$ws.Range("G45").Value = 6100
$ws.Range("H45").Value = 4100
$ws.Range("I45").Value = 3800
$ws.Range("J45").Value = 3700

# --- Cells whose latest-period (J column) figure is now unavailable ---
$ws.Range("J21").Value = "NA"
$ws.Range("J83").Value = "NA"
$ws.Range("J94").Value = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"
